# EDM-4: improved function for updating and reading data from Excel
#
# Rewrites the "Документы" register: row 1 is updated in place, a new
# row 2 is appended, and the date columns (D:E) for a handful of
# additional rows below are pre-formatted as dates (no values yet),
# matching the pattern produced by the document-manager's bulk
# read/update routine.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 3-8: date formatting pre-applied for future entries --------
# Registering this numeric format first mirrors the workbook's existing
# custom-format slot ordering (escaped "dd\.mm\.yyyy" = numFmtId 164).
$ws.Range("D3:E8").NumberFormat = "dd\.mm\.yyyy"

# --- Row 1: existing document, values refreshed ---------------------
$ws.Range("A1").Value = 1
$ws.Range("B1").Value = 312321
$ws.Range("C1").Value = "Доверенность"

$ws.Range("D1").NumberFormat = "dd.mm.yyyy"
$ws.Range("D1").Value = "17.01.2024"

$ws.Range("E1").NumberFormat = "dd.mm.yyyy"
$ws.Range("E1").Value = "26.01.2024"

$ws.Range("F1").Value = 1

# --- Row 2: new document appended ------------------------------------
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 12312
$ws.Range("C2").Value = "fdfd"

# "03.01.2024" reads as an unambiguous day/month/year date once the
# dd.mm.yyyy display format is active, so stage it as Text first and
# switch the format afterwards to keep it a literal date string.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "03.01.2024"
$ws.Range("D2").NumberFormat = "dd.mm.yyyy"

$ws.Range("E2").NumberFormat = "dd.mm.yyyy"
$ws.Range("E2").Value = "26.01.2024"

$ws.Range("F2").Value = 1

# --- Selection / view bookkeeping ------------------------------------
[void]$ws.Range("H10").Select()
